# Auto-generated Excel COM-interop script
# Applies the 'Horarios actualizados Linea 141 - 770' scrape refresh
# to sheets LP1912, LP1912-215 and 6203-6173.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:23:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 43"

$data = @(
  @("00:09:57", "01:11", "215_ALUAR", 62, "LP1912"),
  @("01:17:09", "01:58", "14_ABASTO", 41, "LP1912"),
  @("01:59:40", "02:00", "14_ABASTO", 1, "LP1912"),
  @("02:30:53", "02:57", "215_ALUAR", 27, "LP1912"),
  @("02:55:01", "02:58", "215_ALUAR", 3, "LP1912"),
  @("01:59:40", "03:06", "215_ALUAR", 67, "LP1912"),
  @("03:21:41", "03:24", "14_ABASTO", 3, "LP1912"),
  @("02:55:01", "03:48", "14_ABASTO", 53, "LP1912"),
  @("01:59:40", "03:50", "14_ABASTO", 111, "LP1912"),
  @("02:30:53", "03:52", "14_ABASTO", 82, "LP1912"),
  @("03:57:17", "04:01", "81_EL PELIGRO", 4, "LP1912"),
  @("03:21:41", "04:45", "215A_EL PATO", 84, "LP1912"),
  @("04:46:27", "04:46", "215A_EL PATO", 0, "LP1912"),
  @("04:46:27", "04:53", "11_ETCHEVERRY", 7, "LP1912"),
  @("03:21:41", "05:14", "14_ABASTO", 113, "LP1912"),
  @("04:56:17", "05:16", "17_ROMERO", 20, "LP1912"),
  @("04:56:17", "05:22", "23_HERNANDEZ", 26, "LP1912"),
  @("05:23:05", "05:24", "23_HERNANDEZ", 1, "LP1912"),
  @("04:46:27", "05:34", "215B_EL PATO", 48, "LP1912"),
  @("05:23:05", "05:35", "215B_EL PATO", 12, "LP1912"),
  @("03:57:17", "05:39", "14_ABASTO", 102, "LP1912"),
  @("05:23:05", "05:46", "15_ABASTO", 23, "LP1912"),
  @("04:24:09", "05:48", "215A_EL PATO", 84, "LP1912"),
  @("05:23:05", "05:54", "10_OLMOS", 31, "LP1912"),
  @("05:23:05", "06:04", "16_SANTA ANA", 41, "LP1912"),
  @("04:24:09", "06:09", "16_SANTA ANA", 105, "LP1912"),
  @("04:56:17", "06:11", "215A_EL PATO", 75, "LP1912"),
  @("05:23:05", "06:12", "215A_EL PATO", 49, "LP1912"),
  @("05:23:05", "06:14", "225_HARAS DEL SUR", 51, "LP1912"),
  @("05:23:05", "06:21", "26_HERNANDEZ", 58, "LP1912"),
  @("05:23:05", "06:27", "23_HERNANDEZ", 64, "LP1912"),
  @("04:56:17", "06:29", "86_EST CHICA-ESC AGRARIA", 93, "LP1912"),
  @("05:23:05", "06:30", "86_EST CHICA-ESC AGRARIA", 67, "LP1912"),
  @("05:23:05", "06:31", "16_SANTA ANA", 68, "LP1912"),
  @("05:23:05", "06:44", "225_C ROCA-H SUR", 81, "LP1912"),
  @("04:56:17", "06:46", "215C_EL PATO", 110, "LP1912"),
  @("05:23:05", "06:47", "215C_EL PATO", 84, "LP1912"),
  @("05:23:05", "07:00", "14_ABASTO", 97, "LP1912"),
  @("05:23:05", "07:05", "15_ABASTO", 102, "LP1912"),
  @("05:23:05", "07:07", "225_GOMEZ", 104, "LP1912"),
  @("05:23:05", "07:12", "215A_EL PATO", 109, "LP1912"),
  @("05:23:05", "07:16", "11_ETCHEVERRY", 113, "LP1912"),
  @("05:23:05", "07:21", "26_HERNANDEZ", 118, "LP1912")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:23:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 14"

$data = @(
  @("00:09:57", "01:11", "215_ALUAR", 62, "LP1912"),
  @("02:30:53", "02:57", "215_ALUAR", 27, "LP1912"),
  @("02:55:01", "02:58", "215_ALUAR", 3, "LP1912"),
  @("01:59:40", "03:06", "215_ALUAR", 67, "LP1912"),
  @("03:21:41", "04:45", "215A_EL PATO", 84, "LP1912"),
  @("04:46:27", "04:46", "215A_EL PATO", 0, "LP1912"),
  @("04:46:27", "05:34", "215B_EL PATO", 48, "LP1912"),
  @("05:23:05", "05:35", "215B_EL PATO", 12, "LP1912"),
  @("04:24:09", "05:48", "215A_EL PATO", 84, "LP1912"),
  @("04:56:17", "06:11", "215A_EL PATO", 75, "LP1912"),
  @("05:23:05", "06:12", "215A_EL PATO", 49, "LP1912"),
  @("04:56:17", "06:46", "215C_EL PATO", 110, "LP1912"),
  @("05:23:05", "06:47", "215C_EL PATO", 84, "LP1912"),
  @("05:23:05", "07:12", "215A_EL PATO", 109, "LP1912")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:23:05"
$ws.Cells.Item(3, 1).Value = "Total filas: 5"

$data = @(
  @("00:09:57", "00:09", "215A_LA PLATA", 0, "L6173"),
  @("05:23:05", "05:44", "215A_LA PLATA", 21, "L6173"),
  @("05:23:05", "06:09", "215A_LA PLATA", 46, "L6173"),
  @("05:23:05", "06:33", "215C_LA PLATA", 70, "L6203"),
  @("05:23:05", "07:00", "215B_LP-P MOR-1 Y 57", 97, "L6173")
)
$r = 6
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}
